$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark (was on the title paragraph).
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 2. Locate the paragraph ending with "The score is rounded to two decimals."
$targetParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*The score is rounded to two decimals.*") {
        $targetParagraph = $candidate
        break
    }
}

$r = $targetParagraph.Range
# Exclude the trailing paragraph mark from the range.
$r.End = $r.End - 1
$sentenceEnd = $r.End

$newText = "All the scores are absolute and does not depend on the current number of cities returned. It could have been a possibility but I think it would increase score differences and then possibly mislead the algorithm."

# Append a space followed by the new justification text.
$r.InsertAfter(" " + $newText)

# 3. Force a run break between the original sentence and the new space run
#    by briefly bookmarking just the inserted space, then removing that
#    temporary bookmark (the run split survives the removal).
$spaceStart = $sentenceEnd
$spaceEnd = $sentenceEnd + 1
$spaceRange = $d.Range($spaceStart, $spaceEnd)
$d.Bookmarks.Add("TempSplitMarker", $spaceRange)
$d.Bookmarks.Item("TempSplitMarker").Delete()

# 4. Wrap the newly added justification text with the "_GoBack" bookmark.
$newTextRange = $d.Range($spaceEnd, $r.End)
$d.Bookmarks.Add("_GoBack", $newTextRange)
